$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-10 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-11 Tuesday", 2)

$d.Content.Find.Execute("301÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "383÷9=", 2)
$d.Content.Find.Execute("822÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "391÷4=", 2)
$d.Content.Find.Execute("662÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "321÷6=", 2)
$d.Content.Find.Execute("186÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "833÷6=", 2)
$d.Content.Find.Execute("376÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "933÷4=", 2)
$d.Content.Find.Execute("999÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "182÷6=", 2)
$d.Content.Find.Execute("400÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "751÷3=", 2)
$d.Content.Find.Execute("905÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "258÷2=", 2)
$d.Content.Find.Execute("514÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "435÷3=", 2)
$d.Content.Find.Execute("552÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "780÷8=", 2)
$d.Content.Find.Execute("157÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "555÷8=", 2)
$d.Content.Find.Execute("755÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "492÷7=", 2)
$d.Content.Find.Execute("245÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "350÷9=", 2)
$d.Content.Find.Execute("968÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "925÷7=", 2)
$d.Content.Find.Execute("963÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "493÷2=", 2)
$d.Content.Find.Execute("547÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "317÷7=", 2)
$d.Content.Find.Execute("259÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "466÷3=", 2)
$d.Content.Find.Execute("166÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "113÷4=", 2)
$d.Content.Find.Execute("321÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "843÷2=", 2)
$d.Content.Find.Execute("997÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "946÷3=", 2)
$d.Content.Find.Execute("649÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "298÷7=", 2)
$d.Content.Find.Execute("465÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "887÷8=", 2)
$d.Content.Find.Execute("689÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "181÷8=", 2)
$d.Content.Find.Execute("938÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "470÷9=", 2)
$d.Content.Find.Execute("398÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "502÷8=", 2)
